$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: Binary Search stats updated (Success 5 -> 6, Last Update 45831 -> 45835) ---
$ws.Range("E9").Value = 6
$ws.Range("I9").Value = 45835

# --- Row 18: replace the stub row (A18/B18 only) with a full new entry ---
# "Search in Rotated Sorted Array" (LeetCode 33)
$ws.Range("A18").Value = 33
$ws.Range("B18").Value = "Search in Rotated Sorted Array"
$ws.Range("C18").Value = "#array #binary-search #必背 "
$ws.Range("D18").Value = "medium"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 25

# H18/I18 are dates; copy the date format from the row above first so the
# cell keeps the existing date style (s="3") instead of minting a new one,
# then assign the serial date values.
$ws.Range("H17").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("H18").Value = 45835

$ws.Range("I17").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("I18").Value = 45835

$ws.Rows.Item(18).RowHeight = 34

# --- Sheet view: selection moved to J9, top-left scrolled up a couple rows ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("J9").Select()
